# Maine COVID summary workbook update:
#   - "cases_by_race" sheet: renumber the trailing 2020-12-11 block (rows 27-34)
#     to continue the running index from the previous block, and append a new
#     2020-12-12 / 2020-12-11 block (rows 35-43).
#   - "cases_by_ethnicity" sheet: renumber the trailing 2020-12-11 block
#     (rows 14-16) to continue the running index, and append a new
#     2020-12-12 / 2020-12-11 block (rows 17-19).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: cases_by_race
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("cases_by_race")

# Renumber the running index (column A) of the existing last block (rows 27-34)
# so it continues on from row 26 (24) instead of restarting at 0.
$ws1.Cells.Item(27, 1).Value = 25
$ws1.Cells.Item(28, 1).Value = 26
$ws1.Cells.Item(29, 1).Value = 27
$ws1.Cells.Item(30, 1).Value = 28
$ws1.Cells.Item(31, 1).Value = 29
$ws1.Cells.Item(32, 1).Value = 30
$ws1.Cells.Item(33, 1).Value = 31
$ws1.Cells.Item(34, 1).Value = 32

# New data block: rows 35-43, for DATA_REFRESH_DT 2020-12-12 / DATA_AS_OF_DT 2020-12-11.
$race1 = @(
    @(0, "", 1),
    @(1, "American Indian or Alaska Native", 48),
    @(2, "Asian", 226),
    @(3, "Black or African American", 1318),
    @(4, "Native Hawaiian or Other Pacific Islander", 10),
    @(5, "Not disclosed", 1425),
    @(6, "Other Race", 352),
    @(7, "Two or more", 94),
    @(8, "White", 12146)
)

$destRow = 35
foreach ($rec in $race1) {
    # Clone formatting (incl. the bold/bordered style on column A) from the
    # row above so the appended rows match the existing block's look.
    $ws1.Range("A" + ($destRow - 1) + ":E" + ($destRow - 1)).Copy()
    $ws1.Range("A" + $destRow + ":E" + $destRow).PasteSpecial(-4122) | Out-Null

    # C/D hold date-like strings that must stay literal text, not be
    # reinterpreted as date serials.
    $ws1.Range("C" + $destRow + ":D" + $destRow).NumberFormat = "@"

    $ws1.Cells.Item($destRow, 1).Value = $rec[0]
    $ws1.Cells.Item($destRow, 2).Value = $rec[1]
    $ws1.Cells.Item($destRow, 3).Value = "2020-12-12"
    $ws1.Cells.Item($destRow, 4).Value = "2020-12-11"
    $ws1.Cells.Item($destRow, 5).Value = $rec[2]

    $destRow = $destRow + 1
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet 2: cases_by_ethnicity
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("cases_by_ethnicity")

# Renumber the running index (column A) of the existing last block (rows 14-16)
# so it continues on from row 13 (11) instead of restarting at 0.
$ws2.Cells.Item(14, 1).Value = 12
$ws2.Cells.Item(15, 1).Value = 13
$ws2.Cells.Item(16, 1).Value = 14

# New data block: rows 17-19, for DATA_REFRESH_DT 2020-12-12 / DATA_AS_OF_DT 2020-12-11.
$eth1 = @(
    @(0, "Hispanic or Latino", 351),
    @(1, "Not Hispanic or Latino", 12216),
    @(2, "unknown", 3053)
)

$destRow = 17
foreach ($rec in $eth1) {
    $ws2.Range("A" + ($destRow - 1) + ":E" + ($destRow - 1)).Copy()
    $ws2.Range("A" + $destRow + ":E" + $destRow).PasteSpecial(-4122) | Out-Null

    $ws2.Range("C" + $destRow + ":D" + $destRow).NumberFormat = "@"

    $ws2.Cells.Item($destRow, 1).Value = $rec[0]
    $ws2.Cells.Item($destRow, 2).Value = $rec[1]
    $ws2.Cells.Item($destRow, 3).Value = "2020-12-12"
    $ws2.Cells.Item($destRow, 4).Value = "2020-12-11"
    $ws2.Cells.Item($destRow, 5).Value = $rec[2]

    $destRow = $destRow + 1
}

$excel.CutCopyMode = 0

Write-Output "Maine COVID summary: appended 2020-12-12 blocks to cases_by_race and cases_by_ethnicity."
